$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new values would otherwise be auto-parsed as numbers
# (losing significant trailing zeros / exact text form) so they stay text,
# matching the original inlineStr storage.
foreach ($addr in @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '61.906.44'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").Value = '2.982.50'
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '542.91'
$ws.Range("E5").Value = '  -0.66%  '
$ws.Range("D6").Value = '134.89'
$ws.Range("E6").Value = '  -0.75%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '2.976.10'
$ws.Range("E8").Value = '  -1.37%  '
$ws.Range("D9").Value = '0.485'
$ws.Range("E9").Value = '  -2.57%  '
$ws.Range("D10").Value = '6.69'
$ws.Range("E10").Value = '  +8.42%  '
$ws.Range("D11").Value = '0.146'
$ws.Range("E11").Value = '  -2.42%  '
$ws.Range("D12").Value = '0.444'
$ws.Range("E12").Value = '  -1.83%  '
$ws.Range("D13").Value = '0.0000218'
$ws.Range("E13").Value = '  -2.95%  '
$ws.Range("D14").Value = '33.55'
$ws.Range("E14").Value = '  -2.83%  '
$ws.Range("D15").Value = '3.412.46'
$ws.Range("E15").Value = '  -2.84%  '
$ws.Range("D16").Value = '61.708.82'
$ws.Range("E16").Value = '  -0.59%  '
$ws.Range("D17").Value = '0.108'
$ws.Range("E17").Value = '  -2.14%  '
$ws.Range("D18").Value = '2.970.05'
$ws.Range("E18").Value = '  -1.76%  '
$ws.Range("D19").Value = '6.53'
$ws.Range("E19").Value = '  -2.67%  '
$ws.Range("D20").Value = '463.65'
$ws.Range("E20").Value = '  -3.35%  '
$ws.Range("D21").Value = '13.39'
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").Value = '0.650'
$ws.Range("E22").Value = '  -4.25%  '
$ws.Range("D23").Value = '7.08'
$ws.Range("E23").Value = '  -0.39%  '
$ws.Range("D24").Value = '79.06'
$ws.Range("E24").Value = '  -2.23%  '
$ws.Range("D25").Value = '12.51'
$ws.Range("E25").Value = '  +2.61%  '
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("D27").Value = '2.70'
$ws.Range("E27").Value = '  -1.18%  '
$ws.Range("D28").Value = '7.54'
$ws.Range("E28").Value = '  -3.69%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").Value = '1.97'
$ws.Range("E30").Value = '  +2.74%  '
$ws.Range("D31").Value = '25.13'
$ws.Range("E31").Value = '  -2.99%  '
$ws.Range("D32").Value = '1.12'
$ws.Range("E32").Value = '  -3.12%  '
$ws.Range("D33").Value = '2.30'
$ws.Range("E33").Value = '  -1.16%  '
$ws.Range("D34").Value = '5.47'
$ws.Range("E34").Value = '  +0.17%  '
$ws.Range("D35").Value = '54.11'
$ws.Range("E35").Value = '  -2.51%  '
$ws.Range("D36").Value = '5.79'
$ws.Range("E36").Value = '  -3.02%  '
$ws.Range("D37").Value = '446.41'
$ws.Range("E37").Value = '  -3.18%  '
$ws.Range("D38").Value = '0.0800'
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("D39").Value = '0.0385'
$ws.Range("E39").Value = '  -0.78%  '
$ws.Range("D40").Value = '2.928.91'
$ws.Range("E40").Value = '  -9.28%  '
$ws.Range("E41").Value = '  -5.11%  '
$ws.Range("D42").Value = '7.95'
$ws.Range("E42").Value = '  -2.88%  '
$ws.Range("D43").Value = '2.44'
$ws.Range("E43").Value = '  -1.33%  '
$ws.Range("D44").Value = '26.60'
$ws.Range("E44").Value = '  +3.14%  '
$ws.Range("D46").Value = '0.246'
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("D47").Value = '1.99'
$ws.Range("E47").Value = '  -1.32%  '
$ws.Range("D48").Value = '0.108'
$ws.Range("E48").Value = '  -0.84%  '
$ws.Range("D49").Value = '113.70'
$ws.Range("E49").Value = '  -3.81%  '
$ws.Range("D50").Value = '0.0₃0484'
$ws.Range("E50").Value = '  -3.00%  '
$ws.Range("D51").Value = '1.23'
$ws.Range("E51").Value = '  -2.94%  '
